# Luxury_Housing_Insights.pptx — "Add files via upload" edit:
#   1. Refresh the cached "datetimeFigureOut" field text (9/5/2025 -> 9/18/2025)
#      on the slide master and every slide layout's Date placeholder.
#   2. Fix a typo in slide 5's Observation paragraph: "<450%)" -> "<50%)".

$p = $ppt.ActivePresentation

$oldDate = "9/5/2025"
$newDate = "9/18/2025"
$ppPlaceholderDate = 16

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = -1
        try { $phType = $sh.PlaceholderFormat.Type } catch { $phType = -1 }

        if ($phType -eq $ppPlaceholderDate -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a. Slide master's Date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# 1b. Every slide layout's Date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Update-DatePlaceholder $layouts.Item($l).Shapes
}

# 2. Slide 5, "Content Placeholder 2": fix the conversion-rate typo in the
#    Observation paragraph's last run, without disturbing the other runs.
$slide5 = $p.Slides.Item(5)
$contentShape = $slide5.Shapes.Item(2)
$tr5 = $contentShape.TextFrame.TextRange

$oldRun = " Road, Indira Nagar and Koramangala (all above >50.5%) have higher booking conversion rates compared to Sarjapur Road, Yelahanka and Domlur (<450%)"
$newRun = " Road, Indira Nagar and Koramangala (all above >50.5%) have higher booking conversion rates compared to Sarjapur Road, Yelahanka and Domlur (<50%)"

$startIdx = $tr5.Text.IndexOf($oldRun)
if ($startIdx -ge 0) {
    $runRange = $tr5.Characters($startIdx + 1, $oldRun.Length)
    $runRange.Text = $newRun
}
